$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CargaMasivaDetalle")

# A2: was the shared string "    123A" -> replace with the plain number 123
$ws.Range("A2").Value = 123

# R5: typo fix 6547895 -> 6587945
$ws.Range("R5").Value = 6587945

# New data row 8: same request number as row 2/3/4 (123)
$ws.Range("A8").Value = 123

# Grow the "Tabla1" table (NumeroSolicitud) so it covers the new row
$lo = $ws.ListObjects.Item("Tabla1")
$lo.Resize($ws.Range("A1:A8"))

# Reflect the reviewer's on-screen selection/scroll position after the edit
$ws.Range("R2").Select()
$excel.ActiveWindow.ScrollColumn = 11
$excel.ActiveWindow.ScrollRow = 1

$wb.Save()
